$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# New rows of data (2h each) appended after existing row 81 entries
$ws.Range("B82").Value = 2
$ws.Range("C82").Value = "Kaksi uutta hakukenttää Analyzer työkaluun, parametrien selvittelyä ja yhteisen komponentin DropDownMenu luonti"
$ws.Range("D82").Value = "client"

$ws.Range("B83").Value = 2
$ws.Range("C83").Value = "uusien hakukenttien käyttöönotto haussa, datan näyttäminen uusien arvojen avulla, DataView komponentit"
$ws.Range("D83").Value = "client"

# Apply same style (centered, style index 1) to B82:B83 as other cells in column B
$ws.Range("B82:B83").HorizontalAlignment = -4108
$ws.Range("B82:B83").VerticalAlignment = -4108

# Update SUM formula to include new rows
$ws.Range("B86").Formula = "=SUM(B2:B83)"

# Percentage formula recomputes automatically off the new sum
$ws.Range("B88").Formula = "=B86/B87*100"

$excel.Calculate()

# Update selection to reflect new active cell
$ws.Range("B76").Select()
